# New forms for registration
# Applies the registrationMember.xlsx edit: adds a new "ID Type" select_one
# question (with its own choice list), a new "email" / "vulnerability"
# screen, new translations, and updates to the "model" and "choices" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "survey" (sheet1)
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Row 11 (id_type question) becomes a select_one driven by a new choice list
$survey.Range("C11").Value = "select_one"
$survey.Range("H11").Value = "id_types_list"
$survey.Range("K11").Value = $true
$survey.Range("K11").Font.Size = 18

# Row 22 (is_head_of_household) is now required
$survey.Range("K22").Value = $true
$survey.Range("K22").Font.Size = 18

# New screen asking for email + vulnerability + picture (rows 24-28)
$survey.Range("A10").Copy($survey.Range("A24"))
$survey.Range("A24").Value = "begin screen"
$survey.Range("E10").Copy($survey.Range("E24"))
$survey.Range("F10").Copy($survey.Range("F24"))
$survey.Range("H10").Copy($survey.Range("H24"))
$survey.Range("I10").Copy($survey.Range("I24"))
$survey.Range("J10").Copy($survey.Range("J24"))
$survey.Range("K10").Copy($survey.Range("K24"))

$survey.Range("C6").Copy($survey.Range("C25"))
$survey.Range("D6").Copy($survey.Range("D25"))
$survey.Range("E6").Copy($survey.Range("E25"))
$survey.Range("C25").Value = "string"
$survey.Range("D25").Value = "email"
$survey.Range("E25").Value = "What is your email address?"

$survey.Range("C17").Copy($survey.Range("C26"))
$survey.Range("D17").Copy($survey.Range("D26"))
$survey.Range("E17").Copy($survey.Range("E26"))
$survey.Range("H17").Copy($survey.Range("H26"))
$survey.Range("C26").Value = "select_one"
$survey.Range("D26").Value = "vulnerability"
$survey.Range("E26").Value = "Vulnerability Criteria:"
$survey.Range("H26").Value = "vulnerability_list"

$survey.Range("A18").Copy($survey.Range("A27"))
$survey.Range("A27").Value = "end screen"

$survey.Range("C6").Copy($survey.Range("C28"))
$survey.Range("D6").Copy($survey.Range("D28"))
$survey.Range("E6").Copy($survey.Range("E28"))
$survey.Range("C28").Value = "image"
$survey.Range("D28").Value = "picture"
$survey.Range("E28").Value = "Take a picture of the beneficiary."

$survey.Range("D28").Select()

# ---------------------------------------------------------------------------
# Sheet "table_specific_translations" (sheet3)
# ---------------------------------------------------------------------------
$translations = $wb.Worksheets.Item("table_specific_translations")

$translations.Range("A3").Value = "id_type"
$translations.Range("B3").Value = "ID Type"

$translations.Range("A4").Value = "gender"
$translations.Range("B4").Value = "Gender"

$translations.Range("A5").Value = "age"
$translations.Range("B5").Value = "Age"

$translations.Range("A6").Value = "beneficiary_code"
$translations.Range("B6").Value = "Beneficiary Code"

$translations.Range("A7").Value = "name"
$translations.Range("B7").Value = "Name"

$translations.Range("B12").Select()

# ---------------------------------------------------------------------------
# Sheet "model" (sheet4)
# ---------------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")

$model.Range("A2").Copy($model.Range("A14"))
$model.Range("B2").Copy($model.Range("B14"))
$model.Range("A14").Value = "image"
$model.Range("B14").Value = "picture"

$model.Range("A2").Copy($model.Range("A15"))
$model.Range("B2").Copy($model.Range("B15"))
$model.Range("A15").Value = "string"
$model.Range("B15").Value = "email"

$model.Range("A2").Copy($model.Range("A16"))
$model.Range("B2").Copy($model.Range("B16"))
$model.Range("A16").Value = "string "
$model.Range("B16").Value = "vulnerability"

$model.Range("B13").Select()

# ---------------------------------------------------------------------------
# Sheet "choices" (sheet5)
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Columns.Item(2).ColumnWidth = 21.83203125

$choices.Range("A13").Value = "id_types_list"
$choices.Range("B13").Value = "Registration Number"
$choices.Range("C13").Value = "Registration Number"

$choices.Range("A14").Value = "id_types_list"
$choices.Range("B14").Value = "Case Number"
$choices.Range("C14").Value = "Case Number"

$choices.Range("A15").Value = "id_types_list"
$choices.Range("B15").Value = "Asylum Card"
$choices.Range("C15").Value = "Asylum Card"

$choices.Range("A16").Value = "id_types_list"
$choices.Range("B16").Value = "Police Notes"
$choices.Range("C16").Value = "Police Notes"

$choices.Range("A18").Value = "vulnerability_list"
$choices.Range("B18").Value = "None"
$choices.Range("C18").Value = "None"

$choices.Range("A19").Value = "vulnerability_list"
$choices.Range("B19").Value = "Disability"
$choices.Range("C19").Value = "Disability"

$choices.Range("A20").Value = "vulnerability_list"
$choices.Range("B20").Value = "Unaccompanied minors"
$choices.Range("C20").Value = "Unaccompanied minors"

$choices.Range("A21").Value = "vulnerability_list"
$choices.Range("B21").Value = "Single-parent family"
$choices.Range("C21").Value = "Single-parent family"

$choices.Range("B26").Select()

# Restore "survey" as the active sheet/selection, matching the target workbook
$survey.Activate()
$survey.Range("D28").Select()
